$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("C27").Value = "dc 39.a "
$ws.Range("D27").Value = "coluna  do conj transversal traseiro le"
$ws.Range("E27").Value = "'20"
$ws.Range("F27").Value = "QUALIDADE"
$ws.Range("G27").Value = "izaac"
$ws.Range("H27").Value = "ADM"
$ws.Range("I27").Value = "CMM GLOBAL"
$ws.Range("J27").Value = "INSP DISPOSITIVO"
$ws.Range("K27").Value = "certificação de dispositivo"
$ws.Range("L27").Value = "C2025.0025"

# Row 28
$ws.Range("C28").Value = "dc 39.a "
$ws.Range("D28").Value = "coluna  do conj transversal traseiro le"
$ws.Range("E28").Value = "'20"
$ws.Range("F28").Value = "QUALIDADE"
$ws.Range("G28").Value = "izaac"
$ws.Range("H28").Value = "ADM"
$ws.Range("I28").Value = "CMM GLOBAL"
$ws.Range("J28").Value = "INSP DISPOSITIVO"
$ws.Range("K28").Value = "cert"
$ws.Range("L28").Value = "C2025.0026"
